$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", styled like the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2: numeric value 1
$ws.Range("H2").Value = 1
